$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "before" sheet has a one-row offset bug: several header rows (column A)
# are paired with the WRONG value cells in B/C (stale/misaligned content left
# over from an earlier edit), and some headers (Objectives:, Docentes
# responsaveis:, Syllabus:, Avaliacao:, Requisitos:) have no value row of their
# own at all. Fix this by inserting a new row at 13 - giving "Docentes
# responsaveis:" its own value row - which shifts the old rows 13-23 down to
# 14-24, and then filling in all of the correct content.
$ws.Rows(13).Insert()

# The row-insert leaves a stray empty styled cell at A13 (row 13 only has B/C
# content in the final layout, no header in column A) - remove it.
$ws.Range("A13").Clear()

$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "1. Compreensão dos fenômenos envolvidos no processo de solidificação. 2. Noção dos aspectos práticos do processamento de materiais em empresas de fundição."
$ws.Range("C10").Value = "1. Compreensão dos fenômenos envolvidos no processo de solidificação. 2. Noção dos aspectos práticos do processamento de materiais em empresas de fundição."
$ws.Rows(10).RowHeight = 60

$ws.Range("A11").Value = "Objectives:"
$ws.Rows(11).RowHeight = 60

$ws.Range("A12").Value = "Docentes responsáveis:"

$ws.Range("B13").Value = "5009972 - Gilberto Carvalho Coelho"
$ws.Range("C13").Value = "5009972 - Gilberto Carvalho Coelho"

$ws.Range("A14").Value = "Programa resumido:"
$ws.Range("B14").Value = "A. Introdução; B. Nucleação; C. Redistribuição de soluto na solidificação de ligas; D. Crescimento; E. Macroestruturas de solidificação; F. Aspectos práticos do processamento de materiais por fundição."
$ws.Range("C14").Value = "A. Introdução; B. Nucleação; C. Redistribuição de soluto na solidificação de ligas; D. Crescimento; E. Macroestruturas de solidificação; F. Aspectos práticos do processamento de materiais por fundição."
$ws.Rows(14).RowHeight = 60

$ws.Range("A15").Value = "Short syllabus:"
$ws.Rows(15).RowHeight = 60

$ws.Range("A16").Value = "Programa:"
$ws.Range("B16").Value = "1. Introdução: história da fundição; aplicações e mercado de fundidos; 2. Nucleação: Nucleação homogênea; nucleação heterogênea; taxa de nucleação, agentes nucleantes; 3. Redistribuição de soluto na solidificação: Materiais puros; ligas binárias; coeficiente de redistribuição; solidificação em condições de equilíbrio; solidificação fora de equilíbrio; 4. Crescimento: solidificação de ligas monofásicas - crescimento planar, celular e dendrítico; solidificação de ligas polifásicas ligas eutéticas e ligas peritéticas; 5. Macroestruturas de solidificação: contração volumétrica na solidificação; zonas coquilhada, colunar e equiaxial; controle da macroestrutura; 6. Aspectos práticos do processamento de materiais por fundição: equipamentos e processos de fundição; segregação macro e microssegregação; defeitos originados na solidificação."
$ws.Range("C16").Value = "1. Introdução: história da fundição; aplicações e mercado de fundidos; 2. Nucleação: Nucleação homogênea; nucleação heterogênea; taxa de nucleação, agentes nucleantes; 3. Redistribuição de soluto na solidificação: Materiais puros; ligas binárias; coeficiente de redistribuição; solidificação em condições de equilíbrio; solidificação fora de equilíbrio; 4. Crescimento: solidificação de ligas monofásicas - crescimento planar, celular e dendrítico; solidificação de ligas polifásicas ligas eutéticas e ligas peritéticas; 5. Macroestruturas de solidificação: contração volumétrica na solidificação; zonas coquilhada, colunar e equiaxial; controle da macroestrutura; 6. Aspectos práticos do processamento de materiais por fundição: equipamentos e processos de fundição; segregação macro e microssegregação; defeitos originados na solidificação."
$ws.Rows(16).RowHeight = 120

$ws.Range("A17").Value = "Syllabus:"
$ws.Rows(17).RowHeight = 120

$ws.Range("A18").Value = "Avaliação:"

$ws.Range("A19").Value = "Método:"
$ws.Range("B19").Value = "O curso será ministrado na forma de aulas expositivas. Estão previstas visitas a empresas de fundição para consolidação dos conceitos apresentados nas aulas expositivas."
$ws.Range("C19").Value = "O curso será ministrado na forma de aulas expositivas. Estão previstas visitas a empresas de fundição para consolidação dos conceitos apresentados nas aulas expositivas."
$ws.Rows(19).RowHeight = 60

$ws.Range("A20").Value = "Critério:"
$ws.Range("B20").Value = "Serão aplicadas duas avaliações escritas (P1 e P2) que comporão a nota final (NF). O critério para a nota final é: NF=(P1+P2)/2"
$ws.Range("C20").Value = "Serão aplicadas duas avaliações escritas (P1 e P2) que comporão a nota final (NF). O critério para a nota final é: NF=(P1+P2)/2"
$ws.Rows(20).RowHeight = 60

$ws.Range("A21").Value = "Norma de recuperação:"
$ws.Range("B21").Value = "Será aplicada uma prova de recuperação cuja nota comporá média aritmética com a nota final NF."
$ws.Range("C21").Value = "Será aplicada uma prova de recuperação cuja nota comporá média aritmética com a nota final NF."
$ws.Rows(21).RowHeight = 60

$ws.Range("A22").Value = "Bibliografia:"
$ws.Range("B22").Value = "1. Garcia, A. Solidificação: Fundamentos e Aplicações, Editora da Unicamp, 2001. 2. Flemings, M.C. Solidification Processing, McGraw-Hill, 1974. 3. Pfann, W.G. Zone Melting, John Wiley, 1966. 4. Shewmon, P.G. Diffusion in Solids, McGraw-Hill, 1963. 5. Shewmon, P.G. Transformations in Metals, McGraw-Hill, 1969. 6. Prates, M.; Davis, G.J. Solidificação e Fundição de Metais e suas Ligas, EDUSP, 1978. 7. Davis, G.J. Solidification and Casting, Applied Science Publisher, 1973. 8. Brice, J.C. The Growth of Crystals from the Melt, John Wiley, 1965. 9. Winegard, W.C. An Introduction to Solidification of Metals, Institute of Metals, 1964. 10. Chalmers, B. Principles of Solidification, Robert E. Krieger, 1964. 11. Casting, ASM Handbook, Vol 15, Ninth Edition, ASM International, 1988. 12. Metallography and Microstructures, ASM Handbook, Vol 9, Ninth Edition, ASM International, 1988. 13. Welding, Brasing, and Soldering, ASM Handbook, Vol 6, Ninth Edition, ASM International, 1988."
$ws.Range("C22").Value = "1. Garcia, A. Solidificação: Fundamentos e Aplicações, Editora da Unicamp, 2001. 2. Flemings, M.C. Solidification Processing, McGraw-Hill, 1974. 3. Pfann, W.G. Zone Melting, John Wiley, 1966. 4. Shewmon, P.G. Diffusion in Solids, McGraw-Hill, 1963. 5. Shewmon, P.G. Transformations in Metals, McGraw-Hill, 1969. 6. Prates, M.; Davis, G.J. Solidificação e Fundição de Metais e suas Ligas, EDUSP, 1978. 7. Davis, G.J. Solidification and Casting, Applied Science Publisher, 1973. 8. Brice, J.C. The Growth of Crystals from the Melt, John Wiley, 1965. 9. Winegard, W.C. An Introduction to Solidification of Metals, Institute of Metals, 1964. 10. Chalmers, B. Principles of Solidification, Robert E. Krieger, 1964. 11. Casting, ASM Handbook, Vol 15, Ninth Edition, ASM International, 1988. 12. Metallography and Microstructures, ASM Handbook, Vol 9, Ninth Edition, ASM International, 1988. 13. Welding, Brasing, and Soldering, ASM Handbook, Vol 6, Ninth Edition, ASM International, 1988."
$ws.Rows(22).RowHeight = 120

$ws.Range("A23").Value = "Requisitos:"

$ws.Range("B24").Value = "LOM3005 -  Diagrama de Fases  (Requisito fraco)`n"
$ws.Range("C24").Value = "LOM3005 -  Diagrama de Fases  (Requisito fraco)`n"
$ws.Rows(24).RowHeight = 30

# B13 is a brand-new cell (that row did not exist before the insert above), so it
# picked up the ambiguous/overlapping legacy column-1-2 default style (bold,
# no wrap) instead of the intended column-B value style (regular weight, wrapped).
# Copy the correct formatting over from another column-B value cell so it matches
# exactly (same style, no new style entries created).
$ws.Range("B10").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
